$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 0.1687306666666667
$ws.Range("N2").Value = 0.506192
$ws.Range("O2").Value = 0.05883914643530498
$ws.Range("P2").Value = 0.05883914643530497
$ws.Range("Q2").Value = 0.07658268757688889
$ws.Range("R2").Value = 0.689244188192
$ws.Range("S2").Value = 0.05883914643530498
$ws.Range("T2").Value = 0.05883914643530497

# Row 3 updates
$ws.Range("O3").Value = 0.3874259849494012
$ws.Range("P3").Value = 0.3874259849494012
$ws.Range("S3").Value = 0.3874259849494012
$ws.Range("T3").Value = 0.3874259849494012

# Row 4 updates
$ws.Range("M4").Value = 1.587923333333333
$ws.Range("N4").Value = 4.76377
$ws.Range("O4").Value = 0.5537348686152938
$ws.Range("P4").Value = 0.5537348686152938
$ws.Range("Q4").Value = 0.7207192322244445
$ws.Range("R4").Value = 6.48647309002
$ws.Range("S4").Value = 0.5537348686152938
$ws.Range("T4").Value = 0.5537348686152938
